$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Number formats used throughout the sheet (mirrors xl/styles.xml cellXfs):
#   hh:mm;@  -> start/end time columns (D, E)
#   "0"      -> integer minute totals (F, except the two "sum [h]" rows)
#   "0.00"   -> hour totals (G, and the "sum [h]" / "sum [working weeks]" rows)
$timeFmt = "hh:mm;@"
$intFmt  = "0"
$hrFmt   = "0.00"

# --- Row 75: end time moved 20 minutes later (22:00 -> 22:20) -----------
$ws.Range("E75").Value = 0.93055555555555547

# --- Row 76: was a blank spacer row, now becomes a real data row --------
$ws.Range("A76").Value = 2014
$ws.Range("B76").Value = 3
$ws.Range("C76").Value = 18
$ws.Range("D76").Value = 0.375
$ws.Range("D76").NumberFormat = $timeFmt
$ws.Range("E76").Value = 0.5
$ws.Range("E76").NumberFormat = $timeFmt

$ws.Range("F76").Formula = "=(E76-D76)*24*60"
$ws.Range("F76").NumberFormat = $intFmt
$ws.Range("G76").Formula = "=F76/60"
$ws.Range("G76").NumberFormat = $hrFmt

# --- Row 77: becomes the new blank spacer row (shifted down by one) -----
$ws.Range("D77").Value = $null
$ws.Range("D77").NumberFormat = $timeFmt

# E77 previously held the right-aligned "sum [min]" label (style s=4) -
# clear that formatting entirely before giving it the plain time format
# used by the rest of the D:E columns, so it collapses back to style s=1
# instead of minting a brand-new right-aligned + hh:mm combo style.
$ws.Range("E77").ClearContents() | Out-Null
$ws.Range("E77").ClearFormats() | Out-Null
$ws.Range("E77").NumberFormat = $timeFmt

$ws.Range("F77").ClearContents() | Out-Null
$ws.Range("F77").NumberFormat = $intFmt

# --- Row 78: "sum [min]" summary (was row 77, now shifted to 78) --------
# E78 already carries the right-aligned label style (s=4) from the old
# "sum [h]" label it's replacing, so only the text/format need updating.
$ws.Range("E78").Value = "sum [min]"
$ws.Range("F78").Formula = "=SUM(F2:F77)"
$ws.Range("F78").NumberFormat = $intFmt

# --- Row 79: "sum [h]" summary (was row 78, now shifted to 79) ----------
$ws.Range("E79").Value = "sum [h]"
$ws.Range("F79").Formula = "=F78/60"
$ws.Range("F79").NumberFormat = $hrFmt

# --- Row 80: "sum [working weeks]" summary (new row, mirrors old row 79) -
$ws.Range("E80").Value = "sum [working weeks]"
$ws.Range("E80").HorizontalAlignment = -4152
$ws.Range("F80").Formula = "=F79/38.5"
$ws.Range("F80").NumberFormat = $hrFmt

# --- View state: scroll position + active selection ----------------------
$excel.ActiveWindow.ScrollRow = 64
$ws.Range("E77").Select() | Out-Null
